$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin/Link/Price/Volume table updated by the scheduled GitHub Actions refresh.
# Price-column ("D") values are numeric-looking text (e.g. "67.547.46",
# "0.0000220") that must stay text, so each is written with a temporary
# Text number format to stop Excel re-interpreting it as a number; the
# format/style is restored to Normal immediately after so the cell ends up
# indistinguishable from the untouched cells around it.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.547.46'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.523.40'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '610.99'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.42%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '151.47'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.81%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.522.67'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.12%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -1.15%  '
$ws.Range('E10').Value = '  -1.30%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.05'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.16%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.425'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.93%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000220'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.31%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.117.05'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.19%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '31.95'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.47%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.528.03'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.18%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '67.491.32'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '15.24'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.66%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '447.02'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.75%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.29'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -4.56%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.624'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.54%  '
$ws.Range('E24').Value = '  -0.66%  '
$ws.Range('E25').Value = '  +10.89%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.666.18'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.02%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.18'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -4.82%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.37'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.51'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.77%  '
$ws.Range('E31').Value = '  -4.14%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('E33').Value = '  +4.65%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '25.81'
$ws.Range('D34').Style = "Normal"
$ws.Range('E35').Value = '  -0.51%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.516.61'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.27%  '
$ws.Range('E37').Value = '  -3.70%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '8.06'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.74%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '177.37'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.70%  '
$ws.Range('E42').Value = '  +3.88%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0876'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.34%  '
$ws.Range('E44').Value = '  -3.35%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.881'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.76%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '45.54'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.03%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '27.18'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -5.34%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.61'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.15%  '
$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.25'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +4.03%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.59'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.98%  '
$ws.Range('E51').Value = '  -0.91%  '
